# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "65.794.15"
$ws.Cells.Item(2, 5).Value = "  +1.83%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.481.66"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "581.15"
$ws.Cells.Item(5, 5).Value = "  +0.41%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "161.09"
$ws.Cells.Item(6, 5).Value = "  +1.98%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
$ws.Cells.Item(8, 2).Value = "XRP"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.600"
$ws.Cells.Item(8, 5).Value = "  +7.36%  "

# Row 9
$ws.Cells.Item(9, 2).Value = "LidoStakedEther"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(9, 4).Value = "3.482.70"
$ws.Cells.Item(9, 5).Value = "  +0.45%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "7.28"
$ws.Cells.Item(10, 5).Value = "  -3.96%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.10%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.444"
$ws.Cells.Item(12, 5).Value = "  -1.15%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "4.082.04"
$ws.Cells.Item(13, 5).Value = "  +0.39%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -1.27%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -1.12%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "28.96"
$ws.Cells.Item(16, 5).Value = "  +3.73%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "65.734.05"
$ws.Cells.Item(17, 5).Value = "  +1.68%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.500.15"
$ws.Cells.Item(18, 5).Value = "  +0.90%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.47"
$ws.Cells.Item(19, 5).Value = "  +0.34%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.33"
$ws.Cells.Item(20, 5).Value = "  -0.69%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "391.40"
$ws.Cells.Item(21, 5).Value = "  -1.48%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "8.24"
$ws.Cells.Item(22, 5).Value = "  -3.67%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.41%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "73.50"
$ws.Cells.Item(24, 5).Value = "  +0.31%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.17%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +1.41%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.78"
$ws.Cells.Item(27, 5).Value = "  +1.52%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.63%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 5).Value = "  -0.12%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "6.41"
$ws.Cells.Item(30, 5).Value = "  +5.68%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +2.94%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.65%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "23.77"
$ws.Cells.Item(33, 5).Value = "  -0.44%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.50"
$ws.Cells.Item(34, 5).Value = "  -3.29%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.02%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "7.10"
$ws.Cells.Item(36, 5).Value = "  +0.65%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +4.07%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "163.92"
$ws.Cells.Item(38, 5).Value = "  +2.10%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.96"
$ws.Cells.Item(39, 5).Value = "  +4.29%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "3.080.19"
$ws.Cells.Item(40, 5).Value = "  +5.58%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0770"
$ws.Cells.Item(41, 5).Value = "  -2.73%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "27.28"
$ws.Cells.Item(42, 5).Value = "  -1.50%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0321"
$ws.Cells.Item(43, 5).Value = "  -1.15%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "OKB"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "43.13"
$ws.Cells.Item(44, 5).Value = "  +1.86%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Filecoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "4.52"
$ws.Cells.Item(45, 5).Value = "  +1.62%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.777"
$ws.Cells.Item(46, 5).Value = "  -0.11%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "25.60"
$ws.Cells.Item(47, 5).Value = "  +6.93%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.13"
$ws.Cells.Item(48, 5).Value = "  +2.55%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.26"
$ws.Cells.Item(49, 5).Value = "  +1.76%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Cosmos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "6.70"
$ws.Cells.Item(50, 5).Value = "  +1.83%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Bittensor"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "309.88"
$ws.Cells.Item(51, 5).Value = "  +4.41%  "
